# Add a new data row (row 2) to the "Tài khoản" sheet:
#   A2 = "2", B2 = "a", C2 = "a", F2 = "1"  (all stored as text)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tài khoản")

# A2 -> "2" as text (forces shared-string storage instead of numeric)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2"
$ws.Range("A2").Style = "Normal"

# B2, C2 -> "a"
$ws.Range("B2").Value = "a"
$ws.Range("C2").Value = "a"

# F2 -> "1" as text (forces shared-string storage instead of numeric)
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1"
$ws.Range("F2").Style = "Normal"
